# Append the 2021 data row (row 5) to Sheet1, mirroring the existing
# year rows (2018-2020) already present in rows 2-4.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$targetRow = 5

# Row label, styled the same way as the other year-label cells in
# column A (bold, bordered, centered) by copying A4's formatting.
$ws.Range("A" + $targetRow).Value = "2021年"
$ws.Range("A4").Copy() | Out-Null
$ws.Range("A" + $targetRow).PasteSpecial(-4122) | Out-Null

# Numeric data for the remaining columns (B:DK). Columns G and AG have
# no reported value for this year (left blank, matching rows 2-4).
$rowValues = [ordered]@{
    "B" = -13.4
    "C" = 2.9
    "D" = 73
    "E" = 17.1
    "F" = -4
    "H" = 16.2
    "I" = 6.3
    "J" = 2.3
    "K" = -16
    "L" = -32
    "M" = -11.3
    "N" = -6.7
    "O" = -2.4
    "P" = -11.3
    "Q" = -18
    "R" = -18.9
    "S" = -32.7
    "T" = -9.6
    "U" = -25.7
    "V" = -48.9
    "W" = -9.300000000000001
    "X" = -1.2
    "Y" = 14.8
    "Z" = 2.1
    "AA" = 6.9
    "AB" = 19.1
    "AC" = 7.4
    "AD" = -22.9
    "AE" = 9.5
    "AF" = -5.7
    "AH" = -59.9
    "AI" = -41.8
    "AJ" = -36.9
    "AK" = 26.6
    "AL" = -8.800000000000001
    "AM" = -15.9
    "AN" = -37.7
    "AO" = -38.6
    "AP" = -40
    "AQ" = -40.7
    "AR" = -40.2
    "AS" = -26.4
    "AT" = -38
    "AU" = -10.8
    "AV" = 0.9
    "AW" = -38.6
    "AX" = -8.800000000000001
    "AY" = -17.9
    "AZ" = 8.9
    "BA" = -13.1
    "BB" = -13.8
    "BC" = -5.6
    "BD" = 232.5
    "BE" = -25
    "BF" = -24.4
    "BG" = -13
    "BH" = -49.5
    "BI" = -28.3
    "BJ" = -4.9
    "BK" = -20
    "BL" = -2.6
    "BM" = -2
    "BN" = -21.5
    "BO" = -13.9
    "BP" = 8.699999999999999
    "BQ" = -12.2
    "BR" = 45.3
    "BS" = 4
    "BT" = -1.7
    "BU" = 10.4
    "BV" = 4.4
    "BW" = 9.699999999999999
    "BX" = 10
    "BY" = 0.1
    "BZ" = -24.7
    "CA" = -20.5
    "CB" = 16.2
    "CC" = 31
    "CD" = -9.699999999999999
    "CE" = -28.2
    "CF" = 7.6
    "CG" = -2.9
    "CH" = -80.8
    "CI" = 8.9
    "CJ" = 2.1
    "CK" = 13.6
    "CL" = -26.7
    "CM" = -0.1
    "CN" = 16.7
    "CO" = -6.2
    "CP" = 10.2
    "CQ" = 7.7
    "CR" = -9.300000000000001
    "CS" = -29.9
    "CT" = -12.3
    "CU" = 1.3
    "CV" = -1.8
    "CW" = 78.40000000000001
    "CX" = 8.4
    "CY" = 51.2
    "CZ" = 10.2
    "DA" = 0.7
    "DB" = -8.9
    "DC" = 26.4
    "DD" = -13.3
    "DE" = -21.2
    "DF" = 0.2
    "DG" = 119.4
    "DH" = -4.9
    "DI" = -31.5
    "DJ" = -1.3
    "DK" = 83
}

foreach ($col in $rowValues.Keys) {
    $ws.Range($col + $targetRow).Value = $rowValues[$col]
}

